# Update "想去人数" (want-to-go count) figures in the F column on the
# "展览" and "全部类型" sheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value  = 4578
$wsExhibit.Range("F15").Value = 984
$wsExhibit.Range("F22").Value = 3518
$wsExhibit.Range("F23").Value = 5876
$wsExhibit.Range("F27").Value = 522
$wsExhibit.Range("F32").Value = 2474
$wsExhibit.Range("F40").Value = 1013
$wsExhibit.Range("F41").Value = 909

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 4578
$wsAll.Range("F16").Value = 984
$wsAll.Range("F23").Value = 3518
$wsAll.Range("F24").Value = 5876
$wsAll.Range("F28").Value = 522
$wsAll.Range("F33").Value = 2474
$wsAll.Range("F41").Value = 1013
$wsAll.Range("F42").Value = 909
